# "Add quaters to FOTEST and remove time_quater example"
#
# The Data sheet's "time" column (C) currently holds literal years
# (2018/2019/2020/2021) repeating for every group of 4 rows. Replace them
# with quarter labels 2000Q1..2000Q4 (same cyclic pattern), and make the
# Data sheet the active/selected tab (it was previously Codelists).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$quarters = @("2000Q1", "2000Q2", "2000Q3", "2000Q4")

for ($r = 2; $r -le 85; $r++) {
    $idx = ($r - 2) % 4
    $ws.Cells.Item($r, 3).Value = $quarters[$idx]
}

# Data becomes the active sheet/tab (was Codelists), with a new selection.
$ws.Activate()
$ws.Range("E16").Select()
